$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# buy_orders: swap "fees"/"currency" columns and insert a new "fee_currency"
# column before "comment"
#   OLD: date, symbol, quantity, buy_price, fees, currency, comment
#   NEW: date, symbol, quantity, buy_price, currency, fees, fee_currency, comment
# ---------------------------------------------------------------------------
$wsBuy = $wb.Worksheets.Item("buy_orders")
$wsBuy.Range("E1").Value = "currency"
$wsBuy.Range("F1").Value = "fees"
$wsBuy.Range("G1").Value = "fee_currency"
$wsBuy.Range("G1").Font.Bold = $true
$wsBuy.Range("H1").Value = "comment"
$wsBuy.Range("H1").Font.Bold = $true

# ---------------------------------------------------------------------------
# sell_orders: same header restructuring as buy_orders, plus update the
# existing data rows (2-7): old fees/currency columns swap content and a new
# fee_currency data column is filled in with the (former) currency value.
#   OLD: date, symbol, quantity, sell_price, fees, currency, comment
#   NEW: date, symbol, quantity, sell_price, currency, fees, fee_currency, comment
# ---------------------------------------------------------------------------
$wsSell = $wb.Worksheets.Item("sell_orders")
$wsSell.Range("E1").Value = "currency"
$wsSell.Range("F1").Value = "fees"
$wsSell.Range("G1").Value = "fee_currency"
$wsSell.Range("G1").Font.Bold = $true
$wsSell.Range("H1").Value = "comment"
$wsSell.Range("H1").Font.Bold = $true

for ($r = 2; $r -le 7; $r++) {
    $oldFees = $wsSell.Cells.Item($r, 5).Value2
    $oldCurrency = $wsSell.Cells.Item($r, 6).Value2
    $wsSell.Cells.Item($r, 5).Value = $oldCurrency
    $wsSell.Cells.Item($r, 6).Value = $oldFees
    $wsSell.Cells.Item($r, 7).Value = $oldCurrency
}

# ---------------------------------------------------------------------------
# currency_conversions: rename/reorder columns and append "fees"/"fee_currency"
# before "comment"
#   OLD: date, foreign_amount, source_fees, source_currency, target_currency, comment
#   NEW: date, source_amount, source_currency, target_amount, target_currency, fees, fee_currency, comment
# ---------------------------------------------------------------------------
$wsConv = $wb.Worksheets.Item("currency_conversions")
$wsConv.Range("B1").Value = "source_amount"
$wsConv.Range("C1").Value = "source_currency"
$wsConv.Range("D1").Value = "target_amount"
$wsConv.Range("E1").Value = "target_currency"
$wsConv.Range("F1").Value = "fees"
$wsConv.Range("G1").Value = "fee_currency"
$wsConv.Range("G1").Font.Bold = $true
$wsConv.Range("H1").Value = "comment"
$wsConv.Range("H1").Font.Bold = $true

# ---------------------------------------------------------------------------
# money_transfers: swap "fees"/"currency" and insert "fee_currency" before
# "comment"
#   OLD: date, buy_date, amount, fees, currency, comment
#   NEW: date, buy_date, amount, currency, fees, fee_currency, comment
# ---------------------------------------------------------------------------
$wsTransfer = $wb.Worksheets.Item("money_transfers")
$wsTransfer.Range("D1").Value = "currency"
$wsTransfer.Range("E1").Value = "fees"
$wsTransfer.Range("F1").Value = "fee_currency"
$wsTransfer.Range("F1").Font.Bold = $true
$wsTransfer.Range("G1").Value = "comment"
$wsTransfer.Range("G1").Font.Bold = $true

# ---------------------------------------------------------------------------
# make "sell_orders" the active tab/sheet (previously "money_transfers")
# ---------------------------------------------------------------------------
$wsSell.Activate()
